$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "start/end" entry is being inserted right after row 134 (before the
# blank spacer row that precedes the summary block). Insert a whole row at
# 135 so everything below (spacer + the three summary rows) shifts down by
# one, picking up formatting from the row above as Excel normally does.
$ws.Rows("135:135").Insert()

# Populate the new entry (2014-07-12, 21:00 -> 21:00, i.e. a zero-length
# session) using the same layout as the surrounding rows.
$ws.Range("A135").Value = 2014
$ws.Range("B135").Value = 7
$ws.Range("C135").Value = 12
$ws.Range("D135").Value = 0.875
$ws.Range("E135").Value = 0.875
$ws.Range("F135").Formula = "=(E135-D135)*24*60"
$ws.Range("G135").Formula = "=F135/60"

# The running totals below now need to include the newly inserted row.
$ws.Range("F137").Formula = "=SUM(F2:F135)"
$ws.Range("F138").Formula = "=F137/60"
$ws.Range("F139").Formula = "=F138/38.5"

# Match the recorded post-edit selection.
$ws.Range("E136").Select()
